# case with 380 kV done
# Update the line-loss (pl_mw) results table on Sheet1 with the recomputed
# per-line values for the 380 kV case. Columns A, F, I:N are unchanged
# (index column / zero-valued lines); only B, C, D, E, G, H, O move.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.201376172747416
$ws.Range("C2").Value = 0.3355094820660725
$ws.Range("D2").Value = 0.07970775826393606
$ws.Range("E2").Value = 0.4254940551921464
$ws.Range("G2").Value = 0.2506091721017754
$ws.Range("H2").Value = 0.4034740199273301
$ws.Range("O2").Value = 1.225594442798524

$ws.Range("B3").Value = 1.051949474620756
$ws.Range("C3").Value = 0.2962493627776723
$ws.Range("D3").Value = 0.07213751377541655
$ws.Range("E3").Value = 0.3710793600973403
$ws.Range("G3").Value = 0.2479892266328392
$ws.Range("H3").Value = 0.4076013966952061
$ws.Range("O3").Value = 1.228330647571568

$ws.Range("B4").Value = 0.9599315408196958
$ws.Range("C4").Value = 0.2720377919384589
$ws.Range("D4").Value = 0.06752476887481862
$ws.Range("E4").Value = 0.3377655909772699
$ws.Range("G4").Value = 0.2467861169645147
$ws.Range("H4").Value = 0.4104988740165396
$ws.Range("O4").Value = 1.231618975261696

$ws.Range("B5").Value = 0.9223675857677449
$ws.Range("C5").Value = 0.2621453829302141
$ws.Range("D5").Value = 0.06565390212512057
$ws.Range("E5").Value = 0.3242123819504599
$ws.Range("G5").Value = 0.2463969163119586
$ws.Range("H5").Value = 0.4117706509739065
$ws.Range("O5").Value = 1.233360841055358

$ws.Range("B6").Value = 0.9161261939611904
$ws.Range("C6").Value = 0.2605012030996363
$ws.Range("D6").Value = 0.06534378031780363
$ws.Range("E6").Value = 0.3219631783760804
$ws.Range("G6").Value = 0.2463383681929656
$ws.Range("H6").Value = 0.4119873175448348
$ws.Range("O6").Value = 1.233674271760677

$ws.Range("B7").Value = 0.9594252045629332
$ws.Range("C7").Value = 0.2719044839318485
$ws.Range("D7").Value = 0.06749950187440845
$ws.Range("E7").Value = 0.3375827194749235
$ws.Range("G7").Value = 0.2467804600110384
$ws.Range("H7").Value = 0.410515657478264
$ws.Range("O7").Value = 1.231640843033688

$ws.Range("B8").Value = 1.149910789012324
$ws.Range("C8").Value = 0.32199486484663
$ws.Range("D8").Value = 0.07709015998119639
$ws.Range("E8").Value = 0.4067105053716489
$ws.Range("G8").Value = 0.2496210993426757
$ws.Range("H8").Value = 0.4048215590885036
$ws.Range("O8").Value = 1.22620253512099

$ws.Range("B9").Value = 1.521254469736334
$ws.Range("C9").Value = 0.4193627679769065
$ws.Range("D9").Value = 0.09618131359884785
$ws.Range("E9").Value = 0.5431362864581928
$ws.Range("G9").Value = 0.2584507481109313
$ws.Range("H9").Value = 0.3965518986282035
$ws.Range("O9").Value = 1.228416799267791

$ws.Range("B10").Value = 1.79268744326356
$ws.Range("C10").Value = 0.4903550495919831
$ws.Range("D10").Value = 0.1103861466626199
$ws.Range("E10").Value = 0.6440390492278993
$ws.Range("G10").Value = 0.2669829008021765
$ws.Range("H10").Value = 0.3922618740131583
$ws.Range("O10").Value = 1.238060374397492

$ws.Range("B11").Value = 1.915857821068471
$ws.Range("C11").Value = 0.5225295212978267
$ws.Range("D11").Value = 0.1168883419067157
$ws.Range("E11").Value = 0.6901195367793775
$ws.Range("G11").Value = 0.2713210248802653
$ws.Range("H11").Value = 0.3907023085443342
$ws.Range("O11").Value = 1.244224079252234

$ws.Range("B12").Value = 1.962453946475989
$ws.Range("C12").Value = 0.5346954147681799
$ws.Range("D12").Value = 0.1193564267227885
$ws.Range("E12").Value = 0.7075972006134066
$ws.Range("G12").Value = 0.2730304477858141
$ws.Range("H12").Value = 0.390168461369953
$ws.Range("O12").Value = 1.246816480869057

$ws.Range("B13").Value = 1.952420704130418
$ws.Range("C13").Value = 0.532076076386204
$ws.Range("D13").Value = 0.118824619868235
$ws.Range("E13").Value = 0.7038317931467333
$ws.Range("G13").Value = 0.272659312187443
$ws.Range("H13").Value = 0.39028090625483
$ws.Range("O13").Value = 1.246246625749222

$ws.Range("B14").Value = 1.919692241149733
$ws.Range("C14").Value = 0.5235307783456733
$ws.Range("D14").Value = 0.1170912753469509
$ws.Range("E14").Value = 0.6915568611879905
$ws.Range("G14").Value = 0.2714603183481614
$ws.Range("H14").Value = 0.3906572496726568
$ws.Range("O14").Value = 1.244432163879708

$ws.Range("B15").Value = 1.899639090900166
$ws.Range("C15").Value = 0.5182941907727923
$ws.Range("D15").Value = 0.1160303142020354
$ws.Range("E15").Value = 0.6840418208766152
$ws.Range("G15").Value = 0.2707346114925002
$ws.Range("H15").Value = 0.3908951690158915
$ws.Range("O15").Value = 1.243354481520385

$ws.Range("B16").Value = 1.784631656607644
$ws.Range("C16").Value = 0.4882499016848669
$ws.Range("D16").Value = 0.1099620284264518
$ws.Range("E16").Value = 0.6410313765560716
$ws.Range("G16").Value = 0.2667086701349035
$ws.Range("H16").Value = 0.3923717137446943
$ws.Range("O16").Value = 1.237693558010392

$ws.Range("B17").Value = 1.713998663405562
$ws.Range("C17").Value = 0.4697874784318969
$ws.Range("D17").Value = 0.1062496955009777
$ws.Range("E17").Value = 0.6146933444763505
$ws.Range("G17").Value = 0.2643565820233107
$ws.Range("H17").Value = 0.3933781786772386
$ws.Range("O17").Value = 1.234677862691115

$ws.Range("B18").Value = 1.67334375429698
$ws.Range("C18").Value = 0.4591570944835439
$ws.Range("D18").Value = 0.1041182546929491
$ws.Range("E18").Value = 0.5995611499510005
$ws.Range("G18").Value = 0.2630466610432904
$ws.Range("H18").Value = 0.3939939449407461
$ws.Range("O18").Value = 1.233110384695692

$ws.Range("B19").Value = 1.659573831667672
$ws.Range("C19").Value = 0.4555559087128813
$ws.Range("D19").Value = 0.1033972359356881
$ws.Range("E19").Value = 0.5944404583141818
$ws.Range("G19").Value = 0.2626104890502177
$ws.Range("H19").Value = 0.3942087550858702
$ws.Range("O19").Value = 1.232608271613458

$ws.Range("B20").Value = 1.721520647798172
$ws.Range("C20").Value = 0.4717540072728639
$ws.Range("D20").Value = 0.1066444866088716
$ws.Range("E20").Value = 0.6174953240443273
$ws.Range("G20").Value = 0.2646025148136602
$ws.Range("H20").Value = 0.3932672198373126
$ws.Range("O20").Value = 1.234981576968607

$ws.Range("B21").Value = 1.9293066422224
$ws.Range("C21").Value = 0.5260412294085768
$ws.Range("D21").Value = 0.1176002416839594
$ws.Range("E21").Value = 0.6951615291457358
$ws.Range("G21").Value = 0.2718106746817313
$ws.Range("H21").Value = 0.3905451660811821
$ws.Range("O21").Value = 1.244958081366974

$ws.Range("B22").Value = 2.064838877434227
$ws.Range("C22").Value = 0.5614165811499561
$ws.Range("D22").Value = 0.1247945480334067
$ws.Range("E22").Value = 0.7460850311777563
$ws.Range("G22").Value = 0.2769106813025104
$ws.Range("H22").Value = 0.3890969592404616
$ws.Range("O22").Value = 1.252985431721385

$ws.Range("B23").Value = 1.992527887369988
$ws.Range("C23").Value = 0.5425458412244097
$ws.Range("D23").Value = 0.1209516780324122
$ws.Range("E23").Value = 0.7188904469686861
$ws.Range("G23").Value = 0.2741527882477186
$ws.Range("H23").Value = 0.3898395061199409
$ws.Range("O23").Value = 1.248562241473479

$ws.Range("B24").Value = 1.718120105102969
$ws.Range("C24").Value = 0.4708649896635961
$ws.Range("D24").Value = 0.1064659927423151
$ws.Range("E24").Value = 0.616228518459053
$ws.Range("G24").Value = 0.2644911968344275
$ws.Range("H24").Value = 0.3933172686800361
$ws.Range("O24").Value = 1.23484374999731

$ws.Range("B25").Value = 1.421037123746487
$ws.Range("C25").Value = 0.3931162760424627
$ws.Range("D25").Value = 0.09098569162910053
$ws.Range("E25").Value = 0.5061220949675516
$ws.Range("G25").Value = 0.2463969163119586
$ws.Range("H25").Value = 0.3984769193678233
$ws.Range("O25").Value = 1.226421932574908
